# Penalty Reward System (unfinished) - shift forecast weeks forward by one
# and rewrite the MyForecast numbers, then update the Summary sheet totals.
#
# Note: the Week_Start_Date column (and several Summary values) look like
# dates/numbers, so Excel would otherwise auto-convert them on assignment.
# A leading apostrophe forces them to stay literal text, matching the
# original (inlineStr) cell type.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: Week_Start_Date (B) and MyForecast (D) ---

$ws1.Range("B2").Value  = "'2025-01-12"
$ws1.Range("D2").Value  = 3

$ws1.Range("B3").Value  = "'2025-01-19"
$ws1.Range("D3").Value  = 3

$ws1.Range("B4").Value  = "'2025-01-26"
$ws1.Range("D4").Value  = 3

$ws1.Range("B5").Value  = "'2025-02-02"
$ws1.Range("D5").Value  = 2

$ws1.Range("B6").Value  = "'2025-02-09"
$ws1.Range("D6").Value  = 3

$ws1.Range("B7").Value  = "'2025-02-16"

$ws1.Range("B8").Value  = "'2025-02-23"
$ws1.Range("D8").Value  = 3

$ws1.Range("B9").Value  = "'2025-03-02"

$ws1.Range("B10").Value = "'2025-03-09"
$ws1.Range("D10").Value = 3

$ws1.Range("B11").Value = "'2025-03-16"
$ws1.Range("D11").Value = 3

$ws1.Range("B12").Value = "'2025-03-23"
$ws1.Range("D12").Value = 3

$ws1.Range("B13").Value = "'2025-03-30"
$ws1.Range("D13").Value = 3

$ws1.Range("B14").Value = "'2025-04-06"
$ws1.Range("D14").Value = 3

$ws1.Range("B15").Value = "'2025-04-13"
$ws1.Range("D15").Value = 3

$ws1.Range("B16").Value = "'2025-04-20"
$ws1.Range("D16").Value = 3

$ws1.Range("B17").Value = "'2025-04-27"
$ws1.Range("D17").Value = 3

# --- Summary sheet ---

$ws2.Range("B2").Value  = "2024-02-11 to 2025-01-05"
$ws2.Range("B8").Value  = "123 units"
$ws2.Range("B9").Value  = "'45"
$ws2.Range("B10").Value = "'22"
$ws2.Range("B11").Value = "'11"
$ws2.Range("B12").Value = "'3"
$ws2.Range("B13").Value = "'2025-02-23"
$ws2.Range("B14").Value = "'2"
$ws2.Range("B15").Value = "'2025-02-02"
